# Apply the CVD-file update to the SOR workbook.
# This mirrors the author's regenerated "outputs/SORs/..." values after
# switching automate_finance.qmd to pull the updated CVD files dynamically.

$wb = $excel.ActiveWorkbook

# --- "Aichi Japan" -----------------------------------------------------
$ws = $wb.Worksheets.Item("Aichi Japan")
$ws.Range("L4").Value = 0

# --- "Silvestre Terrazas, Chihuahua " -----------------------------------
$ws = $wb.Worksheets.Item("Silvestre Terrazas, Chihuahua ")
$ws.Range("E3").Value = 0.0776

# --- "La Chaux-de-Fonds Switzerland" ------------------------------------
$ws = $wb.Worksheets.Item("La Chaux-de-Fonds Switzerland")
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

# --- "SEEPZ-SEZ Mumbai India" -------------------------------------------
$ws = $wb.Worksheets.Item("SEEPZ-SEZ Mumbai India")
$ws.Range("G4").Value = 0.0045
$ws.Range("K4").Value = 0.0045
$ws.Range("L4").Value = 0.0046
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("L7").ClearContents()

# --- "West Chester Pennsylvania" ----------------------------------------
$ws = $wb.Worksheets.Item("West Chester Pennsylvania")
$ws.Range("L4").Value = 0.011
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776
$ws.Range("L10").Value = 0.0103
